$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "Based on the condition plot ... cars ... cars ..." paragraph
#   Replace both occurrences of "cars" with "vehicles" AND split the
#   surrounding run into five runs (matching how Word naturally keeps
#   a run boundary around text whose character formatting was touched).
# ---------------------------------------------------------------------

$anchor1 = "Based on the condition plot"
$fullText1 = $d.Content.Text
$anchorIdx1 = $fullText1.IndexOf($anchor1)

# Narrow scope to just this sentence so Find doesn't match "cars"
# elsewhere in the document.
$scopeLen = 300
$scope1 = $d.Range($anchorIdx1, $anchorIdx1 + $scopeLen)

# Locate the first "cars"
$firstCars = $scope1.Duplicate
$firstCars.Find.Execute("cars") | Out-Null
$firstStart = $firstCars.Start
$firstEnd = $firstCars.End

# Locate the second "cars" (search after the first one)
$secondCars = $d.Range($firstEnd, $anchorIdx1 + $scopeLen)
$secondCars.Find.Execute("cars") | Out-Null
$secondStart = $secondCars.Start
$secondEnd = $secondCars.End

# Replace the second occurrence first so the first occurrence's
# offsets remain valid.
$secondRange = $d.Range($secondStart, $secondEnd)
$secondRange.Text = "vehicles"

$firstRange = $d.Range($firstStart, $firstEnd)
$firstRange.Text = "vehicles"

# "vehicles" is 4 characters longer than "cars"; the first replacement
# shifted the second word's start by +4.
$vehiclesLen = 8
$firstVehStart = $firstStart
$firstVehEnd = $firstStart + $vehiclesLen
$secondVehStart = $secondStart + 4
$secondVehEnd = $secondVehStart + $vehiclesLen

# Toggling a character-formatting property on just the "vehicles" span
# forces a run boundary around it (the same mechanism real Word uses
# when it needs a run with different rPr), giving the 5-run structure:
#   ["...the "] ["vehicles"] [" ...by "] ["vehicles"] [" ...chart."]
$rv1 = $d.Range($firstVehStart, $firstVehEnd)
$rv1.Bold = 1
$rv1.Bold = 0

$rv2 = $d.Range($secondVehStart, $secondVehEnd)
$rv2.Bold = 1
$rv2.Bold = 0

# ---------------------------------------------------------------------
# Edit 2: "To compensate the less number of new vehicles ..."
#   Merge the two runs "the " + "less" into a single run "the less"
#   (the surrounding gramStart/gramEnd proofErr markers stay in place).
# ---------------------------------------------------------------------

$anchor2 = "To compensate"
$fullText2 = $d.Content.Text
$anchorIdx2 = $fullText2.IndexOf($anchor2)

$scope2 = $d.Range($anchorIdx2, $anchorIdx2 + 80)
$phrase = $scope2.Duplicate
$phrase.Find.Execute("the less") | Out-Null
$phraseStart = $phrase.Start
$phraseEnd = $phrase.End

# First set it to a placeholder so the engine registers a real change
# (setting identical text is treated as a no-op and would not merge
# the two existing runs into one); then set it back to "the less" so
# the final text matches and the two runs collapse into a single run.
$mergeRange = $d.Range($phraseStart, $phraseEnd)
$mergeRange.Text = "TEMP_PLACEHOLDER"

$mergeRange2 = $d.Range($phraseStart, $phraseStart + 16)
$mergeRange2.Text = "the less"
